$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 data: machine_name / list_value_name / eex_value / is_dataset
$ws.Range("A17").Value = "og_group_ref"
$ws.Range("B17").Value = 144795
$ws.Range("D17").Value = $true

# Move the active selection to D18, matching the saved workbook state
$ws.Range("D18").Select()
